$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '60.894.09'
$ws.Range("E2").Value = '  +3.81%  '

$ws.Range("D3").Value = '2.719.80'
$ws.Range("E3").Value = '  +3.36%  '

$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = '0.999'
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = '  -0.06%  '

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '529.74'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.70%  '

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '147.67'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  +1.96%  '

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.996'
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = '  -0.20%  '

$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.579'
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = '  +0.78%  '

$ws.Range("D9").Value = '2.736.16'
$ws.Range("E9").Value = '  +3.66%  '

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '7.13'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +13.89%  '

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '0.106'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.75%  '

$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '0.342'
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = '  +2.52%  '

$ws.Range("E13").Value = '  +3.10%  '

$ws.Range("D14").Value = '3.194.42'
$ws.Range("E14").Value = '  +3.40%  '

$ws.Range("D15").Value = '60.866.83'
$ws.Range("E15").Value = '  +3.78%  '

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '21.52'
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = '  +3.86%  '

$ws.Range("D17").Value = '2.767.51'
$ws.Range("E17").Value = '  +5.17%  '

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.0000139'
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = '  +1.72%  '

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '345.67'
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = '  -0.48%  '

$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '4.52'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +1.81%  '

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '10.59'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +3.87%  '

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '6.45'
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = '  +5.34%  '

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '1.00'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +0.12%  '

$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '63.44'
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = '  +2.99%  '

$ws.Range("E25").Value = '  +5.09%  '

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '0.419'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +1.47%  '

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '0.994'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  -0.10%  '

$ws.Range("D28").Value = '0.0₃0827'
$ws.Range("E28").Value = '  +3.22%  '

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.33'
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = '  +4.59%  '

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '6.77'
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = '  +8.99%  '

$ws.Range("E31").Value = '  -0.14%  '

$ws.Range("E32").Value = '  +2.26%  '

$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '19.09'
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = '  +1.42%  '

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '150.37'
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = '  +0.91%  '

$ws.Range("E35").Value = '  +7.26%  '

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '1.23'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +8.45%  '

$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.927'
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = '  -6.03%  '

$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.906'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  +8.16%  '

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '1.54'
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = '  +8.88%  '

$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '37.35'
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = '  +2.06%  '

$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '3.69'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +2.13%  '

$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.626'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +4.61%  '

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '282.02'
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = '  +0.85%  '

$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '20.28'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  +4.14%  '

$ws.Range("B45").Value = 'Stellar'
$ws.Range("C45").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.0988'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +0.65%  '

$ws.Range("B46").Value = 'FirstDigitalUSD'
$ws.Range("C46").Value = 'https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd'
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '0.995'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  -0.15%  '

$ws.Range("D47").Value = '2.114.57'
$ws.Range("E47").Value = '  +6.44%  '

$ws.Range("B48").Value = 'Hedera'
$ws.Range("C48").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '0.0544'
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = '  +4.29%  '

$ws.Range("B49").Value = 'RenderToken'
$ws.Range("C49").Value = 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr'
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '4.91'
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = '  +5.90%  '

$ws.Range("E50").Value = '  +2.30%  '

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '19.47'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +6.12%  '
